$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "report generated" timestamp in M6 (date serial number)
$ws.Range("M6").Value = 41439.324970023146

# Populate the three new data rows (13, 14, 15) mirroring row 12's layout.
# Values that look like numbers/dates/booleans need a leading apostrophe so
# they are stored as literal text (matching the source data), exactly like
# the sibling cells in row 12 (e.g. A12 = "002201301").
$ws.Range("A13").Value = "'00120130607"
$ws.Range("C13").Value = "GP/DK//00003/2013"
$ws.Range("E13").Value = "'2013-06-03"
$ws.Range("G13").Value = "unknown"
$ws.Range("I13").Value = "'2013-06-03"
$ws.Range("K13").Value = "unknown"
$ws.Range("M13").Value = "'false"
$ws.Range("O13").Value = "'0"

$ws.Range("A14").Value = "'00120130607"
$ws.Range("C14").Value = "GP/DK//00001/2013"
$ws.Range("E14").Value = "'2013-06-03"
$ws.Range("G14").Value = "unknown"
$ws.Range("I14").Value = "'2013-06-03"
$ws.Range("K14").Value = "unknown"
$ws.Range("M14").Value = "'false"
$ws.Range("O14").Value = "'0"

$ws.Range("A15").Value = "'00120130607"
$ws.Range("C15").Value = "GP/DK//00002/2013"
$ws.Range("E15").Value = "'2013-06-03"
$ws.Range("G15").Value = "unknown"
$ws.Range("I15").Value = "'2013-06-03"
$ws.Range("K15").Value = "unknown"
$ws.Range("M15").Value = "'false"
$ws.Range("O15").Value = "'0"
